$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: change Mating_type value from "a" to "alp", update Nb_gen value
$ws.Range("C2").Value = "alp"
$ws.Range("F2").Value = 5.0432241748072615

# Row 3: new row, duplicate of row2 but with Replicate "r2" and a new Nb_gen value
$ws.Range("A3").Value = "CN"
$ws.Range("B3").Value = "F1"
$ws.Range("C3").Value = "alp"
$ws.Range("D3").Value = "r2"
$ws.Range("E3").Value = "T1"
$ws.Range("F3").Value = 5.0433554250353074

# Update selection to match target view state
$ws.Range("C7").Select()
